$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50:51").Insert()

$ws.Range("F50").Value = "Audit"
$ws.Range("E50").Value = "BuildIt Binary"
$ws.Range("E51").Value = "BuildIt Monetary"
$ws.Range("G51").Value = "BuildIt's monetary audit population consisting of 3500 transactions."
$ws.Range("G50").Value = "BuildIt's non-monetary audit population consisting of 3500 records. "

$ws.Range("D50").Value = "Yes"
$ws.Range("H50").Value = "No"
$ws.Range("I50").Value = "Yes"

$ws.Range("D51").Value = "Yes"
$ws.Range("F51").Value = "Audit"
$ws.Range("H51").Value = "No"
$ws.Range("I51").Value = "Yes"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Columns("H").ColumnWidth = 18.1666666666667

$ws.Range("I52").Select()
